$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.006.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").Value = "'3.238.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'597.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").Value = "'136.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.77%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "'3.237.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").Value = "'0.511"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("D10").Value = "'0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").Value = "'5.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.98%  "

$ws.Range("D13").Value = "'0.0000240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.25%  "

$ws.Range("D14").Value = "'33.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.26%  "

$ws.Range("D15").Value = "'3.794.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "'3.258.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "'63.121.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").Value = "'6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").Value = "'467.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").Value = "'13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").Value = "'0.721"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "

$ws.Range("D23").Value = "'7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.79%  "

$ws.Range("D24").Value = "'13.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'84.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'2.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.91%  "

$ws.Range("D30").Value = "'6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.79%  "

$ws.Range("D31").Value = "'2.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").Value = "'27.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").Value = "'0.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.60%  "

$ws.Range("D34").Value = "'2.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.76%  "

$ws.Range("D35").Value = "'1.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("D36").Value = "'5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.92%  "

$ws.Range("D37").Value = "'51.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.31%  "

$ws.Range("D38").Value = "'0.0₃0718"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").Value = "'423.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "

$ws.Range("D41").Value = "'3.038.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("E42").Value = "  +5.98%  "

$ws.Range("D43").Value = "'8.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.85%  "

$ws.Range("D44").Value = "'2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.12%  "

$ws.Range("D45").Value = "'0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.89%  "

$ws.Range("D46").Value = "'2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "

$ws.Range("D47").Value = "'36.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.58%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'126.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.88%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'25.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("E51").Value = "  -1.98%  "
